$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Duplicate Sheet1 -> Sheet2, placed right after Sheet1. The copy becomes
# the active/selected sheet and Sheet1 loses its tabSelected flag, matching
# the diff automatically.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# Fill in the header cells that turn the two plain ranges into tables.
$ws2.Range("B5").Value = "Column1"
$ws2.Range("C5").Value = "Column2"
$ws2.Range("B9").Value = "Column1"
$ws2.Range("C9").Value = "Column2"

# Turn each block into a real Excel Table.
$t1 = $ws2.ListObjects.Add(1, $ws2.Range("A5:C6"), 0, 1)
$t1.Name = "Table1"

$t2 = $ws2.ListObjects.Add(1, $ws2.Range("A9:C10"), 0, 1)
$t2.Name = "Table2"

# Match the new sheet's recorded selection.
$ws2.Range("D11").Select()
